$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts alleinerz..uhv one column right:
# C..Q becomes D..R), to make room for the new "pid" column.
$ws.Columns.Item(3).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 3).Value = "pid"

# pid values for data rows 2-12 (odd numbers 23..43).
$pidValues = @(23, 25, 27, 29, 31, 33, 35, 37, 39, 41, 43)
for ($i = 0; $i -lt $pidValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $pidValues[$i]
}

# Update the active selection to match the post-edit state (G10 instead of Q10).
$ws.Range("G10").Select()

$wb.Save()
